$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "39.525.67"
$ws.Cells.Item(2, 5).Value = "  +0.60%  "
$ws.Cells.Item(3, 4).Value = "2.162.58"
$ws.Cells.Item(3, 5).Value = "  +0.44%  "
$ws.Cells.Item(4, 5).Value = "  +0.18%  "
$ws.Cells.Item(5, 5).Value = "  -0.94%  "
$ws.Cells.Item(6, 5).Value = "  -0.14%  "
$ws.Cells.Item(7, 4).Value = "62.62"
$ws.Cells.Item(7, 5).Value = "  -0.08%  "
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 5).Value = "  -0.92%  "
$ws.Cells.Item(10, 5).Value = "  -0.99%  "
$ws.Cells.Item(11, 5).Value = "  +0.28%  "
$ws.Cells.Item(12, 5).Value = "  -1.12%  "
$ws.Cells.Item(13, 4).Value = "2.482.10"
$ws.Cells.Item(13, 5).Value = "  +0.53%  "
$ws.Cells.Item(14, 4).Value = "21.65"
$ws.Cells.Item(14, 5).Value = "  -2.87%  "
$ws.Cells.Item(15, 5).Value = "  -1.42%  "
$ws.Cells.Item(16, 5).Value = "  -1.99%  "
$ws.Cells.Item(17, 4).Value = "2.161.83"
$ws.Cells.Item(17, 5).Value = "  +1.30%  "
$ws.Cells.Item(18, 4).Value = "39.510.28"
$ws.Cells.Item(18, 5).Value = "  +0.54%  "
$ws.Cells.Item(19, 4).Value = "71.54"
$ws.Cells.Item(19, 5).Value = "  -0.76%  "
$ws.Cells.Item(20, 4).Value = "0.0₃0885"
$ws.Cells.Item(20, 5).Value = "  +3.84%  "
$ws.Cells.Item(21, 4).Value = "5.99"
$ws.Cells.Item(21, 5).Value = "  -2.40%  "
$ws.Cells.Item(22, 4).Value = "227.34"
$ws.Cells.Item(22, 5).Value = "  -0.39%  "
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 4).Value = "2.34"
$ws.Cells.Item(24, 5).Value = "  +0.20%  "
$ws.Cells.Item(25, 5).Value = "  -3.96%  "
$ws.Cells.Item(26, 4).Value = "170.23"
$ws.Cells.Item(27, 5).Value = "  -3.23%  "
$ws.Cells.Item(28, 5).Value = "  -0.49%  "
$ws.Cells.Item(29, 5).Value = "  +2.40%  "
$ws.Cells.Item(30, 4).Value = "19.62"
$ws.Cells.Item(30, 5).Value = "  -0.22%  "
$ws.Cells.Item(31, 4).Value = "2.68"
$ws.Cells.Item(31, 5).Value = "  +4.21%  "
$ws.Cells.Item(32, 5).Value = "  +0.28%  "
$ws.Cells.Item(33, 4).Value = "4.46"
$ws.Cells.Item(33, 5).Value = "  -3.05%  "
$ws.Cells.Item(34, 5).Value = "  -2.89%  "
$ws.Cells.Item(35, 5).Value = "  -2.87%  "
$ws.Cells.Item(36, 4).Value = "0.0616"
$ws.Cells.Item(36, 5).Value = "  -0.43%  "
$ws.Cells.Item(37, 4).Value = "3.80"
$ws.Cells.Item(37, 5).Value = "  +6.78%  "
$ws.Cells.Item(38, 4).Value = "2.39"
$ws.Cells.Item(38, 5).Value = "  -1.03%  "
$ws.Cells.Item(39, 5).Value = "  +0.16%  "
$ws.Cells.Item(40, 4).Value = "4.96"
$ws.Cells.Item(40, 5).Value = "  +19.07%  "
$ws.Cells.Item(41, 4).Value = "101.97"
$ws.Cells.Item(41, 5).Value = "  -1.04%  "
$ws.Cells.Item(42, 4).Value = "0.0226"
$ws.Cells.Item(42, 5).Value = "  -1.70%  "
$ws.Cells.Item(43, 5).Value = "  -1.85%  "
$ws.Cells.Item(44, 4).Value = "1.511.52"
$ws.Cells.Item(44, 5).Value = "  -1.58%  "
$ws.Cells.Item(45, 5).Value = "  +0.99%  "
$ws.Cells.Item(46, 4).Value = "7.87"
$ws.Cells.Item(46, 5).Value = "  +0.26%  "
$ws.Cells.Item(47, 5).Value = "  -0.20%  "
$ws.Cells.Item(48, 4).Value = "0.0914"
$ws.Cells.Item(48, 5).Value = "  -0.65%  "
$ws.Cells.Item(49, 5).Value = "  -1.57%  "
$ws.Cells.Item(50, 4).Value = "0.000196"
$ws.Cells.Item(50, 5).Value = "  +32.22%  "
$ws.Cells.Item(51, 4).Value = "2.98"
$ws.Cells.Item(51, 5).Value = "  +0.32%  "
